$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time")

# Clear the stray empty inline-string cell at D17 so it becomes truly blank.
$ws.Range("D17").Value = $null

# New rows captured by TimeKeeper, dated 16-Jun-2025.
$rows = @(
    @{ Row = 18; A = "16-Jun-2025"; B = "New";  C = "0h:00m"; D = $null;          E = "12:06 PM"; F = "12:06 PM"; G = 3 },
    @{ Row = 19; A = "16-Jun-2025"; B = "New";  C = "0h:00m"; D = "Test notes";   E = "12:07 PM"; F = "12:07 PM"; G = 5 },
    @{ Row = 20; A = "16-Jun-2025"; B = "New";  C = "0h:00m"; D = $null;          E = "12:09 PM"; F = "12:09 PM"; G = 3 },
    @{ Row = 21; A = "16-Jun-2025"; B = "New";  C = "0h:00m"; D = $null;          E = "12:11 PM"; F = "12:11 PM"; G = 4 },
    @{ Row = 22; A = "16-Jun-2025"; B = "New";  C = "0h:00m"; D = "dfghj";        E = "12:11 PM"; F = "12:11 PM"; G = 4 },
    @{ Row = 23; A = "16-Jun-2025"; B = "Test"; C = "0h:00m"; D = "test";         E = "12:12 PM"; F = "12:13 PM"; G = 6 },
    @{ Row = 24; A = "16-Jun-2025"; B = "New";  C = "0h:00m"; D = "test";         E = "12:13 PM"; F = "12:13 PM"; G = 4 }
)

foreach ($r in $rows) {
    # Column A holds dates formatted like "16-Jun-2025" as plain text in this
    # workbook (matching every prior row) -- force text formatting first so
    # the assignment doesn't get auto-parsed into a date serial number.
    $ws.Cells.Item($r.Row, 1).NumberFormat = "@"
    $ws.Cells.Item($r.Row, 1).Value = $r.A

    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    if ($null -ne $r.D) {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}
